# "Create Emp Test Case"
# Adds a new "employee" worksheet (after the existing "login"/"User" sheets)
# populated with a small employee/account fixture, matching the OrangeHRM
# test-data workbook used by the automation suite.

$wb = $excel.ActiveWorkbook

# Update the selection on the previously-active "User" sheet first, so that
# once we create and select the new sheet below, "employee" ends up as the
# workbook's active/selected tab (matches tabSelected moving to the new sheet).
$ws2 = $wb.Worksheets.Item("User")
[void]$ws2.Range("C23").Select()

# Add the new worksheet at the end of the workbook and name it "employee".
$wsEmp = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsEmp.Name = "employee"

# Header/data block for the first two columns.
$wsEmp.Range("A1").Value = "FirstName"
$wsEmp.Range("B1").Value = "LastName"
$wsEmp.Range("A2").Value = "Navdeep"
$wsEmp.Range("B2").Value = "Kaur"

# Remaining header row.
$wsEmp.Range("C1").Value = "UserName"
$wsEmp.Range("D1").Value = "Password"
$wsEmp.Range("E1").Value = "Status"

# Remaining data row.
$wsEmp.Range("C2").Value = "Navdeep"
$wsEmp.Range("D2").Value = "admin123"
$wsEmp.Range("E2").Value = "Enabled"

# Leave the selection on the new sheet where the author left it.
[void]$wsEmp.Range("B1").Select()
